$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header columns I (I0) and J (IF) - copy formatting from H1 (existing header style)
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for columns I and J, rows 2-16
$values = @(
    @(8, 9),
    @(8, 9),
    @(7, 9),
    @(8, 9),
    @(9, 9),
    @(9, 9),
    @(6, 7),
    @(9, 9),
    @(5, 7),
    @(7, 7),
    @(5, 5),
    @(7, 8),
    @(6, 7),
    @(5, 6),
    @(3, 3)
)

for ($i = 0; $i -lt $values.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $values[$i][0]
    $ws.Cells.Item($row, 10).Value = $values[$i][1]
}
